# Adds a new quarterly column (CF) to the real-time GDP series sheet,
# mirroring the pattern already present for column CE:
#   - CF1 gets the next vintage date (2025-11-25, serial 45986) with the
#     same header formatting as CE1.
#   - CF4:CF35 receive a copy of the most-recent (CE) value/format for
#     each row, matching the existing "carry the latest vintage forward"
#     pattern used throughout the sheet. Rows 2, 3 and 36 have no data in
#     column CE, so they stay empty in column CF as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell (row 1): new vintage date, same style as CE1 ---
$ws.Range("CE1").Copy()
$ws.Range("CF1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("CF1").Value2 = 45986        # 2025-11-25

# --- Data rows 4-35: copy value + format from column CE into column CF ---
for ($r = 4; $r -le 35; $r++) {
    $srcCell = $ws.Cells.Item($r, 83)  # column CE
    $dstCell = $ws.Cells.Item($r, 84)  # column CF
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4163)       # xlPasteAll
}
